$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15309.44921875
$ws.Cells.Item(2, 3).Value = 15302.58203125
$ws.Cells.Item(3, 2).Value = 14884.005859375
$ws.Cells.Item(3, 3).Value = 14794.0947265625
$ws.Cells.Item(4, 2).Value = 14768.142578125
$ws.Cells.Item(4, 3).Value = 14540.65234375
$ws.Cells.Item(5, 2).Value = 13538.70703125
$ws.Cells.Item(5, 3).Value = 13830.166015625
$ws.Cells.Item(6, 2).Value = 12859.1103515625
$ws.Cells.Item(6, 3).Value = 13013.0673828125
$ws.Cells.Item(7, 2).Value = 12269.5048828125
$ws.Cells.Item(7, 3).Value = 12371.6123046875
$ws.Cells.Item(8, 2).Value = 11926.3623046875
$ws.Cells.Item(8, 3).Value = 12090.4033203125
$ws.Cells.Item(9, 2).Value = 12001.630859375
$ws.Cells.Item(9, 3).Value = 12138.8408203125
$ws.Cells.Item(10, 2).Value = 12175.1884765625
$ws.Cells.Item(10, 3).Value = 12307.7802734375
$ws.Cells.Item(11, 2).Value = 12608.7919921875
$ws.Cells.Item(11, 3).Value = 12608.69921875
$ws.Cells.Item(12, 2).Value = 13181.7724609375
$ws.Cells.Item(12, 3).Value = 13094.1015625
$ws.Cells.Item(13, 2).Value = 14199.2421875
$ws.Cells.Item(13, 3).Value = 13948.109375
$ws.Cells.Item(14, 2).Value = 14466.3671875
$ws.Cells.Item(14, 3).Value = 14621.91015625
$ws.Cells.Item(15, 2).Value = 14994.5361328125
$ws.Cells.Item(15, 3).Value = 15063.03125
$ws.Cells.Item(16, 2).Value = 15449.3779296875
$ws.Cells.Item(16, 3).Value = 15427.765625
$ws.Cells.Item(17, 2).Value = 16157.73828125
$ws.Cells.Item(17, 3).Value = 15959.0390625
$ws.Cells.Item(18, 2).Value = 16467.64453125
$ws.Cells.Item(18, 3).Value = 16412.861328125
$ws.Cells.Item(19, 2).Value = 17199.10546875
$ws.Cells.Item(19, 3).Value = 16931.021484375
$ws.Cells.Item(20, 2).Value = 17669.318359375
$ws.Cells.Item(20, 3).Value = 17469.56640625
$ws.Cells.Item(21, 2).Value = 17784.96875
$ws.Cells.Item(21, 3).Value = 17735.212890625
$ws.Cells.Item(22, 2).Value = 18067.65234375
$ws.Cells.Item(22, 3).Value = 17861.19140625
$ws.Cells.Item(23, 2).Value = 17953.484375
$ws.Cells.Item(23, 3).Value = 17835.322265625
$ws.Cells.Item(24, 2).Value = 17743.22265625
$ws.Cells.Item(24, 3).Value = 17610.302734375
$ws.Cells.Item(25, 2).Value = 17284.7421875
$ws.Cells.Item(25, 3).Value = 17173.865234375
$ws.Cells.Item(26, 2).Value = 16520.587890625
$ws.Cells.Item(26, 3).Value = 16488.74609375
$ws.Cells.Item(27, 2).Value = 16365.666015625
$ws.Cells.Item(27, 3).Value = 16004.28515625
$ws.Cells.Item(28, 2).Value = 15898.8076171875
$ws.Cells.Item(28, 3).Value = 15656.7021484375
$ws.Cells.Item(29, 2).Value = 14290.2900390625
$ws.Cells.Item(29, 3).Value = 14627.13671875
$ws.Cells.Item(30, 2).Value = 13360.5302734375
$ws.Cells.Item(30, 3).Value = 13508.5732421875
$ws.Cells.Item(31, 2).Value = 12814.4111328125
$ws.Cells.Item(31, 3).Value = 12794.5791015625
$ws.Cells.Item(32, 2).Value = 12220.5966796875
$ws.Cells.Item(32, 3).Value = 12438.80078125
$ws.Cells.Item(33, 2).Value = 12137.9599609375
$ws.Cells.Item(33, 3).Value = 12415.9873046875
$ws.Cells.Item(34, 2).Value = 12481.8505859375
$ws.Cells.Item(34, 3).Value = 12636.30859375
$ws.Cells.Item(35, 2).Value = 13028.03515625
$ws.Cells.Item(35, 3).Value = 13052.232421875
$ws.Cells.Item(36, 2).Value = 13578.73828125
$ws.Cells.Item(36, 3).Value = 13576.318359375
$ws.Cells.Item(37, 2).Value = 14667.8779296875
$ws.Cells.Item(37, 3).Value = 14437.9091796875
$ws.Cells.Item(38, 2).Value = 15072.005859375
$ws.Cells.Item(38, 3).Value = 15187.740234375
$ws.Cells.Item(39, 2).Value = 15405.0126953125
$ws.Cells.Item(39, 3).Value = 15582.5771484375
$ws.Cells.Item(40, 2).Value = 15702.408203125
$ws.Cells.Item(40, 3).Value = 15764.1669921875
$ws.Cells.Item(41, 2).Value = 16175.1669921875
$ws.Cells.Item(41, 3).Value = 16043.0771484375
$ws.Cells.Item(42, 2).Value = 16147.7255859375
$ws.Cells.Item(42, 3).Value = 16182.9326171875
$ws.Cells.Item(43, 2).Value = 16707.806640625
$ws.Cells.Item(43, 3).Value = 16445.6953125
$ws.Cells.Item(44, 2).Value = 16907.0546875
$ws.Cells.Item(44, 3).Value = 16767.56640625
$ws.Cells.Item(45, 2).Value = 16881.302734375
$ws.Cells.Item(45, 3).Value = 16870.74609375
$ws.Cells.Item(46, 2).Value = 16523.4140625
$ws.Cells.Item(46, 3).Value = 16602.150390625
$ws.Cells.Item(47, 2).Value = 16474.478515625
$ws.Cells.Item(47, 3).Value = 16334.7236328125
$ws.Cells.Item(48, 2).Value = 16711.22265625
$ws.Cells.Item(48, 3).Value = 16392.23828125
$ws.Cells.Item(49, 2).Value = 16275.908203125
$ws.Cells.Item(49, 3).Value = 16264.2197265625
$ws.Cells.Item(50, 2).Value = 16309.34375
$ws.Cells.Item(50, 3).Value = 16045.505859375
$ws.Cells.Item(51, 2).Value = 16042.373046875
$ws.Cells.Item(51, 3).Value = 15820.392578125
$ws.Cells.Item(52, 2).Value = 15565.4931640625
$ws.Cells.Item(52, 3).Value = 15421.640625
$ws.Cells.Item(53, 2).Value = 14542.8525390625
$ws.Cells.Item(53, 3).Value = 14595.0224609375
$ws.Cells.Item(54, 2).Value = 13323.4716796875
$ws.Cells.Item(54, 3).Value = 13655.8388671875
$ws.Cells.Item(55, 2).Value = 12919.025390625
$ws.Cells.Item(55, 3).Value = 12838.3134765625
$ws.Cells.Item(56, 2).Value = 12366.013671875
$ws.Cells.Item(56, 3).Value = 12429.8125
$ws.Cells.Item(57, 2).Value = 12220.0439453125
$ws.Cells.Item(57, 3).Value = 12411.171875
$ws.Cells.Item(58, 2).Value = 12272.83203125
$ws.Cells.Item(58, 3).Value = 12497.9267578125
$ws.Cells.Item(59, 2).Value = 12706.513671875
$ws.Cells.Item(59, 3).Value = 12739.5986328125
$ws.Cells.Item(60, 2).Value = 13308.158203125
$ws.Cells.Item(60, 3).Value = 13224.740234375
$ws.Cells.Item(61, 2).Value = 14209.0966796875
$ws.Cells.Item(61, 3).Value = 14032.8544921875
$ws.Cells.Item(62, 2).Value = 14645.0732421875
$ws.Cells.Item(62, 3).Value = 14745.5380859375
$ws.Cells.Item(63, 2).Value = 15061.591796875
$ws.Cells.Item(63, 3).Value = 15191.1611328125
$ws.Cells.Item(64, 2).Value = 15523.5927734375
$ws.Cells.Item(64, 3).Value = 15531.5126953125
$ws.Cells.Item(65, 2).Value = 15813.302734375
$ws.Cells.Item(65, 3).Value = 15811.1455078125
$ws.Cells.Item(66, 2).Value = 16275.3681640625
$ws.Cells.Item(66, 3).Value = 16140.1015625
$ws.Cells.Item(67, 2).Value = 16761.080078125
$ws.Cells.Item(67, 3).Value = 16575.341796875
$ws.Cells.Item(68, 2).Value = 17039.498046875
$ws.Cells.Item(68, 3).Value = 16950.1015625
$ws.Cells.Item(69, 2).Value = 17033.767578125
$ws.Cells.Item(69, 3).Value = 17052.845703125
$ws.Cells.Item(70, 2).Value = 17021.45703125
$ws.Cells.Item(70, 3).Value = 16982.546875
$ws.Cells.Item(71, 2).Value = 16921.580078125
$ws.Cells.Item(71, 3).Value = 16859.794921875
$ws.Cells.Item(72, 2).Value = 16900.0234375
$ws.Cells.Item(72, 3).Value = 16767.224609375
$ws.Cells.Item(73, 2).Value = 16301.7314453125
$ws.Cells.Item(73, 3).Value = 16391.3671875
$ws.Cells.Item(74, 2).Value = 16043.4462890625
$ws.Cells.Item(74, 3).Value = 15936.6806640625
$ws.Cells.Item(75, 2).Value = 15906.072265625
$ws.Cells.Item(75, 3).Value = 15684.7861328125
$ws.Cells.Item(76, 2).Value = 15463.123046875
$ws.Cells.Item(76, 3).Value = 15386.953125
$ws.Cells.Item(77, 2).Value = 14154.21484375
$ws.Cells.Item(77, 3).Value = 14465.091796875
$ws.Cells.Item(78, 2).Value = 13147.7841796875
$ws.Cells.Item(78, 3).Value = 13431.837890625
$ws.Cells.Item(79, 2).Value = 12596.3232421875
$ws.Cells.Item(79, 3).Value = 12672.2685546875
$ws.Cells.Item(80, 2).Value = 11986.1181640625
$ws.Cells.Item(80, 3).Value = 12216.478515625
$ws.Cells.Item(81, 2).Value = 11993.59375
$ws.Cells.Item(81, 3).Value = 12190.916015625
$ws.Cells.Item(82, 2).Value = 12077.7255859375
$ws.Cells.Item(82, 3).Value = 12320.669921875
$ws.Cells.Item(83, 2).Value = 12603.6435546875
$ws.Cells.Item(83, 3).Value = 12637.25390625
$ws.Cells.Item(84, 2).Value = 13373.482421875
$ws.Cells.Item(84, 3).Value = 13255.783203125
$ws.Cells.Item(85, 2).Value = 14331.125
$ws.Cells.Item(85, 3).Value = 14178.8056640625
$ws.Cells.Item(86, 2).Value = 14777.9775390625
$ws.Cells.Item(86, 3).Value = 14928.783203125
$ws.Cells.Item(87, 2).Value = 15116.654296875
$ws.Cells.Item(87, 3).Value = 15315.4833984375
$ws.Cells.Item(88, 2).Value = 15507.1025390625
$ws.Cells.Item(88, 3).Value = 15558.4296875
$ws.Cells.Item(89, 2).Value = 16269.2421875
$ws.Cells.Item(89, 3).Value = 16058.67578125
$ws.Cells.Item(90, 2).Value = 16776.046875
$ws.Cells.Item(90, 3).Value = 16660.869140625
$ws.Cells.Item(91, 2).Value = 17613.888671875
$ws.Cells.Item(91, 3).Value = 17372.96875
$ws.Cells.Item(92, 2).Value = 18133.810546875
$ws.Cells.Item(92, 3).Value = 18018.6796875
$ws.Cells.Item(93, 2).Value = 18474.994140625
$ws.Cells.Item(93, 3).Value = 18429.6953125
$ws.Cells.Item(94, 2).Value = 18649.736328125
$ws.Cells.Item(94, 3).Value = 18601.75390625
$ws.Cells.Item(95, 2).Value = 18502.5234375
$ws.Cells.Item(95, 3).Value = 18512.81640625
$ws.Cells.Item(96, 2).Value = 17970.1796875
$ws.Cells.Item(96, 3).Value = 18069.619140625
$ws.Cells.Item(97, 2).Value = 17430.943359375
$ws.Cells.Item(97, 3).Value = 17451.5546875
$ws.Cells.Item(98, 2).Value = 16901.16796875
$ws.Cells.Item(98, 3).Value = 16870.142578125
$ws.Cells.Item(99, 2).Value = 16484.77734375
$ws.Cells.Item(99, 3).Value = 16403.740234375
$ws.Cells.Item(100, 2).Value = 16145.9931640625
$ws.Cells.Item(100, 3).Value = 16031.2509765625
$ws.Cells.Item(101, 2).Value = 14651.4775390625
$ws.Cells.Item(101, 3).Value = 15071.205078125
